$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$abstractText = @"
Background
id="Par1">The average length of stay (LOS) in the intensive care unit (ICU_ALOS) is a helpful parameter summarizing critical bed occupancy.

 During the outbreak of a novel virus, estimating early a reliable ICU_ALOS estimate of infected patients is critical to accurately parameterize models examining mitigation and preparedness scenarios.


Methods
id="Par2">Two estimation methods of ICU_ALOS were compared: the average LOS of already discharged patients at the date of estimation (DPE), and a standard parametric method used for analyzing time-to-event data which fits a given distribution to observed data and includes the censored stays of patients still treated in the ICU at the date of estimation (CPE).

 Methods were compared on a series of all COVID-19 consecutive cases (n$([char]0x2009)=$([char]0x2009)59) admitted in an ICU devoted to such patients.

 At the last follow-up date, 99$([char]0x00A0)days after the first admission, all patients but one had been discharged.

 A simulation study investigated the generalizability of the methods' patterns.

 CPE and DPE estimates were also compared to COVID-19 estimates reported to date.


Results
id="Par3">LOS$([char]0x2009)$([char]0x2265)$([char]0x2009)30$([char]0x00A0)days concerned 14 out of the 59 patients (24%), including 8 of the 21 deaths observed.

 Two months after the first admission, 38 (64%) patients had been discharged, with corresponding DPE and CPE estimates of ICU_ALOS (95% CI) at 13.0$([char]0x00A0)days (10.4$([char]0x2013)15.6) and 23.1$([char]0x00A0)days (18.1$([char]0x2013)29.7), respectively.

 Series' true ICU_ALOS was greater than 21$([char]0x00A0)days, well above reported estimates to date.


Conclusions
id="Par4">Discharges of short stays are more likely observed earlier during the course of an outbreak.

 Cautious unbiased ICU_ALOS estimates suggest parameterizing a higher burden of ICU bed occupancy than that adopted to date in COVID-19 forecasting models.


Funding
id="Par5">Support by the National Natural Science Foundation of China (81900097 to Dr.

 Zhou) and the Emergency Response Project of Hubei Science and Technology Department (2020FCA023 to Pr.

 Zhao).



"@

$authorsText = @"
[Nathanael%Lapidus%nathanael.lapidus@inserm.fr%1,    Xianlong%Zhou%xianlongzhou@whu.edu.cn%2,    Xianlong%Zhou%xianlongzhou@whu.edu.cn%0,    Fabrice%Carrat%fabrice.carrat@iplesp.upmc.fr%2,    Fabrice%Carrat%fabrice.carrat@iplesp.upmc.fr%0,    Bruno%Riou%bruno.riou@aphp.fr%2,    Bruno%Riou%bruno.riou@aphp.fr%0,    Yan%Zhao%doctoryanzhao@whu.edu.cn%3,    Yan%Zhao%doctoryanzhao@whu.edu.cn%0,    Gilles%Hejblum%gilles.hejblum@inserm.fr%2,    Gilles%Hejblum%gilles.hejblum@inserm.fr%0]
"@

$ws.Range("D2").Value = $abstractText
$ws.Range("E2").Value = $authorsText

# Re-run autofit on row 2 so the simulated "grow row to fit newlines" side
# effect of writing multi-line content through the object model doesn't
# leave a spurious custom row-height override (the source workbook never
# hard-codes per-row heights even for its longest abstracts).
$ws.Rows.Item(2).AutoFit() | Out-Null
